$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for every existing data row
# (rows 2-310) from 45192 (2023-09-23) to 45202 (2023-10-03).
$ws.Range("C2:C310").Value = 45202

# Row 310 picks up an explicit row height (15pt, custom) in the new file.
$ws.Rows.Item(310).RowHeight = 15

# Append the new record as row 311.
$ws.Range("A311").Value = "A 46828-2023"

$ws.Range("B311").Value = 45199
$ws.Range("B311").NumberFormat = "YYYY-MM-DD"

$ws.Range("C311").Value = 45202
$ws.Range("C311").NumberFormat = "YYYY-MM-DD"

$ws.Range("D311").Value = "HALLANDS LÄN"
$ws.Range("E311").Value = "LAHOLM"

$ws.Range("G311").Value = 9.5
$ws.Range("H311:Q311").Value = 0

$ws.Range("R311").WrapText = $true
